# Mark the three "Lista de adyacencia" / "Matriz de adyacencia" rows as
# done on the "Entregables" sheet (column C uses the shared "x" marker).
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Entregables")

$ws1.Range("C16").Value = "x"
$ws1.Range("C17").Value = "x"
$ws1.Range("C18").Value = "x"

# Leave the selection on the last cell touched.
$ws1.Activate()
$ws1.Range("C18").Select()

# On "Rúbrica", the cursor ends up back at B1.
$ws2 = $wb.Worksheets.Item("Rúbrica")
$ws2.Activate()
$ws2.Range("B1").Select()

# Finish up on the "Método" sheet, which becomes the active tab.
$ws3 = $wb.Worksheets.Item("Método")
$ws3.Activate()
$ws3.Range("E10").Select()
